$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: PAYOUT% 1% -> 0%, TC PAYOUT 2315.44 -> 0
$ws.Range("C2").Value = "'0%"
$ws.Range("D2").Value = 0

# Row 3: PAYOUT% 3% -> 1%, TC PAYOUT 9993.12 -> 1096.55
$ws.Range("C3").Value = "'1%"
$ws.Range("D3").Value = 1096.55

# Row 4: PAID FEEDBACK MORE THAN DEMAND -> LESS THAN DEMAND, PAYOUT% 0% -> 2%, TC PAYOUT 0 -> 4803.74
$ws.Range("B4").Value = "LESS THAN DEMAND"
$ws.Range("C4").Value = "'2%"
$ws.Range("D4").Value = 4803.74

# Row 5: PAYOUT% 2% -> 0%, TC PAYOUT 2810.98 -> 0
$ws.Range("C5").Value = "'0%"
$ws.Range("D5").Value = 0
